$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.007.05"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "2.251.34"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "270.33"
$ws.Range("E5").Value = "  +4.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.06"
$ws.Range("E6").Value = "  +11.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.19"
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0923"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.64"
$ws.Range("E12").Value = "  +8.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "2.595.03"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.91"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "2.249.83"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "43.915.52"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.02"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.53"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.38"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.76"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("E24").Value = "  -3.39%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("E26").Value = "  +12.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.56"
$ws.Range("E28").Value = "  +6.75%  "
$ws.Range("E29").Value = "  +5.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.77"
$ws.Range("E30").Value = "  -6.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.72"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0901"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.91"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("E36").Value = "  +4.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0350"
$ws.Range("E37").Value = "  -4.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.37"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.45"
$ws.Range("E39").Value = "  +19.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.20"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.56"
$ws.Range("E41").Value = "  -4.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.06"
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.204"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.56"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.87"
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0984"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.20"
$ws.Range("E48").Value = "  +4.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.14"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.50"
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.430"
$ws.Range("E51").Value = "  -9.19%  "
